$d = $word.ActiveDocument

# Locate the sentence "... After two epochs of utilizing the widely
# accepted one cycle learning policy ..." and narrow the found range down
# to just the word "two" so only that word gets replaced.
$found = $d.Content
$ok = $found.Find.Execute("two epochs", $true, $false, $false, $false, `
                           $false, $true, 1, $false, $null, 0)

if ($ok) {
    $wordRng = $d.Range($found.Start, $found.Start + 3)
    $wordRng.Text = "five"

    # Word tracks the location of the most recent edit with the built-in
    # "_GoBack" bookmark; move it to sit right after the word we just
    # typed, matching where the insertion point would land after typing.
    $goBackRng = $d.Range($wordRng.End, $wordRng.End)
    $d.Bookmarks.Add("_GoBack", $goBackRng) | Out-Null
}
